# Updated symbol list on Mon Jan 23 08:50:05 UTC 2023 with GitHub Actions
# Refresh coin prices / 1h volume %, and re-sync the GateToken/MXToken/... row
# ordering (rows 8-17) to match the latest coinranking.com snapshot.
# Leading "'" forces text (matches the original inlineStr cell type) so
# numeric-looking values like "304.33" and percentages like "0.70%" are not
# coerced into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.33"
$ws.Range("E2").Value = "'0.70%"
$ws.Range("D3").Value = "'35.86"
$ws.Range("E3").Value = "'-4.24%"
$ws.Range("D4").Value = "'5.089"
$ws.Range("E4").Value = "'1.77%"
$ws.Range("D5").Value = "'0.07851"
$ws.Range("E5").Value = "'0.34%"
$ws.Range("D6").Value = "'2.114"
$ws.Range("E6").Value = "'-3.26%"
$ws.Range("D7").Value = "'7.938"
$ws.Range("E7").Value = "'-1.28%"
$ws.Range("B8").Value = "'MXToken"
$ws.Range("C8").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9193"
$ws.Range("E8").Value = "'1.24%"
$ws.Range("B9").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.09749"
$ws.Range("E9").Value = "'0.82%"
$ws.Range("B10").Value = "'WazirX"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1855"
$ws.Range("E10").Value = "'-1.91%"
$ws.Range("B11").Value = "'MandalaExchangeToken"
$ws.Range("C11").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08629"
$ws.Range("E11").Value = "'1.55%"
$ws.Range("B12").Value = "'BitrueCoin"
$ws.Range("C12").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03559"
$ws.Range("E12").Value = "'1.10%"
$ws.Range("B13").Value = "'BitMartToken"
$ws.Range("C13").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09933"
$ws.Range("E13").Value = "'-0.27%"
$ws.Range("B14").Value = "'BitForexToken"
$ws.Range("C14").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001439"
$ws.Range("E14").Value = "'-2.67%"
$ws.Range("B15").Value = "'TigerCash"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005666"
$ws.Range("E15").Value = "'0.27%"
$ws.Range("B16").Value = "'LEO"
$ws.Range("C16").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.446"
$ws.Range("E16").Value = "'-0.49%"
$ws.Range("B17").Value = "'GateToken"
$ws.Range("C17").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.105"
$ws.Range("E17").Value = "'1.75%"
$ws.Range("D18").Value = "'2.514"
$ws.Range("E18").Value = "'21.29%"
$ws.Range("D19").Value = "'0.3423"
$ws.Range("E19").Value = "'-1.16%"
$ws.Range("D20").Value = "'5.242"
$ws.Range("E20").Value = "'10.05%"
$ws.Range("D21").Value = "'0.1312"
$ws.Range("E21").Value = "'1.45%"
$ws.Range("D22").Value = "'0.2201"
$ws.Range("E22").Value = "'-0.19%"
$ws.Range("D23").Value = "'0.04555"
$ws.Range("E23").Value = "'-1.48%"
$ws.Range("D24").Value = "'0.005069"
$ws.Range("E24").Value = "'5.46%"
$ws.Range("D25").Value = "'0.001233"
$ws.Range("E25").Value = "'0.46%"
$ws.Range("D27").Value = "'0.0004751"
$ws.Range("E27").Value = "'0.07%"
$ws.Range("D39").Value = "'0.01843"
$ws.Range("E39").Value = "'4.91%"
$ws.Range("D40").Value = "'0.04722"
$ws.Range("E40").Value = "'0.02%"
$ws.Range("D41").Value = "'0.007554"
$ws.Range("E41").Value = "'-6.26%"
$ws.Range("D42").Value = "'0.1396"
$ws.Range("E42").Value = "'0.14%"
$ws.Range("D43").Value = "'0.007731"
$ws.Range("E43").Value = "'0.94%"
$ws.Range("D44").Value = "'0.002205"
$ws.Range("E44").Value = "'1.60%"
$ws.Range("D45").Value = "'0.01100"
$ws.Range("E45").Value = "'5.46%"
$ws.Range("D46").Value = "'0.00006318"
$ws.Range("E46").Value = "'4.09%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.09%"
$ws.Range("D48").Value = "'0.0005812"
$ws.Range("E48").Value = "'0.19%"
$ws.Range("E49").Value = "'462.12%"
$ws.Range("D50").Value = "'0.002001"
$ws.Range("E50").Value = "'-25.64%"
$ws.Range("D51").Value = "'0.00002101"
$ws.Range("E51").Value = "'0.09%"
